$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append a second learning-log entry (row 3): branch, date, hours, summary,
# and reference links - mirroring the structure of the existing row 2.
# ---------------------------------------------------------------------------

# Force the date-like text to stay literal text (no auto date conversion),
# same as the existing "08/03/2024" entry in B2.
$ws.Range("B3").NumberFormat = "@"

$ws.Range("A3").Value = "test-2"
$ws.Range("B3").Value = "08/06/2024"
$ws.Range("C3").Value = "9-4:30"
$ws.Range("D3").Value = "Suspense - loading UI - streaming with Suspense - error - global-error - redirect function - permanentRedirect function - redirect with useRouter() - redirect in next.config.js - route groups - project organization - dynamic routes - catch-all segments - parallel routes - slots - tab groups - modals - intercepting routes"

$refs2 = "https://nextjs.org/docs/app/building-your-application/routing/loading-ui-and-streaming - https://nextjs.org/docs/app/building-your-application/routing/error-handling - https://nextjs.org/docs/app/building-your-application/routing/redirecting - https://nextjs.org/docs/app/building-your-application/routing/route-groups - https://nextjs.org/docs/app/building-your-application/routing/colocation - https://nextjs.org/docs/app/building-your-application/routing/dynamic-routes - https://nextjs.org/docs/app/building-your-application/routing/parallel-routes - https://nextjs.org/docs/app/building-your-application/routing/intercepting-routes"
$ws.Range("E3").Value = $refs2

# Copy the row-2 cell formatting down onto the new row so fonts/alignment/
# borders all match the existing table styling.
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)

# "What I Learned"/"References" columns are left-aligned (in addition to the
# existing top-aligned, wrapped text) for both data rows.
$ws.Range("D2:E3").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Turn the References cells into live hyperlinks (first URL of the list is
# the link target; the full multi-link text remains the visible text).
# ---------------------------------------------------------------------------

$refs1 = "https://nextjs.org/docs/getting-started/installation`nhttps://nextjs.org/docs/getting-started/project-structure`nhttps://nextjs.org/docs/app/building-your-application/routing/defining-routes`nhttps://nextjs.org/docs/app/building-your-application/routing/pages-and-layouts`nhttps://nextjs.org/docs/app/building-your-application/routing/linking-and-navigating"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://nextjs.org/docs/getting-started/installation", "", "", $refs1)
$ws.Hyperlinks.Add($ws.Range("E3"), "https://nextjs.org/docs/app/building-your-application/routing/loading-ui-and-streaming", "", "", $refs2)

# Adding a hyperlink re-styles the cell with the blue/underlined Hyperlink
# font; restore the plain "What I Learned" column look used elsewhere.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Resize rows/columns so the longer content is fully visible.
# ---------------------------------------------------------------------------

$ws.Rows.Item(1).RowHeight = 22.5
$ws.Rows.Item(2).RowHeight = 165.75
$ws.Rows.Item(3).RowHeight = 306

$ws.Columns.Item(1).ColumnWidth = 24.75
$ws.Columns.Item(2).ColumnWidth = 19.33
$ws.Columns.Item(3).ColumnWidth = 19.75
$ws.Columns.Item(4).ColumnWidth = 31.75
$ws.Columns.Item(5).ColumnWidth = 35.5

# Leave the selection on the newly added row, matching the saved view state.
$ws.Range("A3").Select()
